$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Structure YCOL_SETTING_S" -----------------------------------
# A new row describing the "SUM_TOTAL" / "Checkbox - Sum column Yes/No"
# field was inserted right above the existing "ALIGNMENT" row (old row 9).
$ws2 = $wb.Worksheets.Item(2)

# Insert a blank row at row 9; this shifts the old rows 9-14 down to 10-15
# while keeping every existing shared-string reference intact.
$ws2.Rows.Item(9).Insert()
# The insert also pushed the trailing blank row (old row 14) into row 15;
# remove it again so the sheet still ends at row 14, exactly like before.
$ws2.Rows.Item(15).Delete()

# Copy the formatting (borders/font) of the row below onto the new row so
# it matches the rest of the table (thin border box, same font/style).
$ws2.Range("A10:E10").Copy()
$ws2.Range("A9:E9").PasteSpecial(-4122)

# Fill in the new row's values (string cell filled first so the new shared
# strings are appended in the same order as the reference workbook).
$ws2.Range("E9").Value = "Checkbox - Sum column Yes/No"
$ws2.Range("A9").Value = "SUM_TOTAL"
$ws2.Range("B9").Value = "XFELD"
$ws2.Range("C9").Value = "XFELD"
$ws2.Range("D9").Value = 1

# --- Selection / active-sheet state -----------------------------------------
# The saved workbook now opens on the "Structure YCOL_SETTING_S" tab with
# E10 selected (previously "Methods Description" / H7 was active).
$ws2.Range("E10").Select()
